$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Extensoes para as Ferramentas" activity (old row 6) was dropped from
# the schedule, so remove its whole row; everything below moves up one row.
$ws.Rows.Item(6).Delete()

# A couple of the neighbouring activities' highlighted (in-progress) date
# spans were touched up after the shift. Reproduce that by copying cell
# formatting (fill/border/font) from an already-correctly-styled donor cell
# onto the cells that need to change, so the saved file reuses a matching
# style instead of inventing a new one.

# Row 6 ("Finalizacao do texto...") / Row 7 ("Defesa...") on column F drop
# back to the plain (non-highlighted) look used elsewhere in the sheet.
$ws.Range("F3").Copy()
$ws.Range("F6:F7").PasteSpecial(-4122)

# Row 5 ("Pesquisa com Profissionais...") span shifts right: E5:G5 drop out
# of the highlight ...
$ws.Range("F3").Copy()
$ws.Range("E5:G5").PasteSpecial(-4122)

# ... Row 4 ("Caracterizacao...") grows into G4:H4 ...
$ws.Range("E3").Copy()
$ws.Range("G4:H4").PasteSpecial(-4122)

# ... and I5:J5 pick the highlight up.
$ws.Range("E3").Copy()
$ws.Range("I5:J5").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Move the active selection to where the user last clicked.
$ws.Range("F7").Select()
